$wb = $excel.ActiveWorkbook

# ---- Metadata sheet: bump the "Date" value ----
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-09-12T13:34:32+00:00"

# ---- Elements sheet: append "|4.0.1" FHIR version pins to type/binding references ----
$els = $wb.Worksheets.Item("Elements")

$els.Range("K10").Value = "canonical(StructureDefinition|4.0.1)`n"
$els.Range("Z11").Value = "http://hl7.org/fhir/ValueSet/security-labels|4.0.1"
$els.Range("Z12").Value = "http://hl7.org/fhir/ValueSet/common-tags|4.0.1"
$els.Range("Z14").Value = "http://hl7.org/fhir/ValueSet/languages|4.0.1"
$els.Range("K20").Value = "Reference(CarePlan|4.0.1|DeviceRequest|4.0.1|ImmunizationRecommendation|4.0.1|MedicationRequest|4.0.1|NutritionOrder|4.0.1|ServiceRequest|4.0.1)`n"
$els.Range("K21").Value = "Reference(MedicationAdministration|4.0.1|MedicationDispense|4.0.1|MedicationStatement|4.0.1|Procedure|4.0.1|Immunization|4.0.1|ImagingStudy|4.0.1)`n"
$els.Range("Z23").Value = "http://hl7.org/fhir/ValueSet/observation-category|4.0.1"
$els.Range("Z24").Value = "http://hl7.org/fhir/ValueSet/observation-codes|4.0.1"
$els.Range("K25").Value = "Reference(Patient|4.0.1|Group|4.0.1|Device|4.0.1|Location|4.0.1)`n"
$els.Range("K26").Value = "Reference(Resource|4.0.1)`n"
$els.Range("K27").Value = "Reference(Encounter|4.0.1)`n"
$els.Range("K30").Value = "Reference(Practitioner|4.0.1|PractitionerRole|4.0.1|Organization|4.0.1|CareTeam|4.0.1|Patient|4.0.1|RelatedPerson|4.0.1)`n"
$els.Range("Z32").Value = "http://hl7.org/fhir/ValueSet/data-absent-reason|4.0.1"
$els.Range("Z33").Value = "http://hl7.org/fhir/ValueSet/observation-interpretation|4.0.1"
$els.Range("Z35").Value = "http://hl7.org/fhir/ValueSet/body-site|4.0.1"
$els.Range("Z36").Value = "http://hl7.org/fhir/ValueSet/observation-methods|4.0.1"
$els.Range("K37").Value = "Reference(Specimen|4.0.1)`n"
$els.Range("K38").Value = "Reference(Device|4.0.1|DeviceMetric|4.0.1)`n"
$els.Range("K43").Value = "Quantity {SimpleQuantity|4.0.1}`n"
$els.Range("K44").Value = "Quantity {SimpleQuantity|4.0.1}`n"
$els.Range("Z45").Value = "http://hl7.org/fhir/ValueSet/referencerange-meaning|4.0.1"
$els.Range("Z46").Value = "http://hl7.org/fhir/ValueSet/referencerange-appliesto|4.0.1"
$els.Range("K49").Value = "Reference(Observation|4.0.1|QuestionnaireResponse|4.0.1|MolecularSequence|4.0.1)`n"
$els.Range("K50").Value = "Reference(DocumentReference|4.0.1|ImagingStudy|4.0.1|Media|4.0.1|QuestionnaireResponse|4.0.1|Observation|4.0.1|MolecularSequence|4.0.1)`n"
$els.Range("Z55").Value = "http://hl7.org/fhir/ValueSet/observation-codes|4.0.1"
$els.Range("Z57").Value = "http://hl7.org/fhir/ValueSet/data-absent-reason|4.0.1"
$els.Range("Z58").Value = "http://hl7.org/fhir/ValueSet/observation-interpretation|4.0.1"

# Column widths K (11) and Z (26) grow to fit the now-longer text (bestFit autofit).
# (requested widths are pre-compensated for this host's char-width/pixel rounding
# so the stored OOXML width lands as close as possible to the authoring tool's value)
$els.Columns.Item(11).ColumnWidth = 121.66666666666667
$els.Columns.Item(26).ColumnWidth = 46.833333333333336

# Re-assert the columns that were already hidden (C, D, AE, AF, AG) so a full-sheet
# re-serialize of <cols> doesn't silently reveal them.
$els.Columns.Item(3).Hidden = $true
$els.Columns.Item(4).Hidden = $true
$els.Columns.Item(31).Hidden = $true
$els.Columns.Item(32).Hidden = $true
$els.Columns.Item(33).Hidden = $true
